$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell text (shared string reshuffle in OOXML is an
# internal detail; what matters is the visible value of A1).
$ws.Range("A1").Value = "Item To Search"

# Move the active selection to B12, matching the authored change.
$ws.Range("B12").Select()
